$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# New row 12: date of the entry, problem solved, and status.
# Copy A11's number format (date, dd/mm/yyyy style) onto A12 first so the
# new date cell reuses the existing date style instead of minting a new one.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = Get-Date -Year 2019 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0

$ws.Range("B12").Value = "51. N-Queens"
$ws.Range("E12").Value = "done"

# Match the saved selection state (active cell moved to D13).
$ws.Range("D13").Select()
